$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 16, shifting rows 16:40 down to 17:41
$ws.Rows.Item(16).Insert()

# Fill the new row 16 with values (constants copied from the rest of the table, new data values)
$ws.Cells.Item(16, 1).Value = 8
$ws.Cells.Item(16, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(16, 3).Value = "Coquimbo"
$ws.Cells.Item(16, 4).Value = 44469
$ws.Cells.Item(16, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16, 5).Value = 4
$ws.Cells.Item(16, 6).Value = 100112052
$ws.Cells.Item(16, 7).Value = "Albahaca"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 700
$ws.Cells.Item(16, 11).Value = 4000
$ws.Cells.Item(16, 12).Value = 4500
$ws.Cells.Item(16, 13).Value = 4250
$ws.Cells.Item(16, 14).Value = "$/paquete"
$ws.Cells.Item(16, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(16, 16).Value = 4250
$ws.Cells.Item(16, 17).Value = 1
$ws.Cells.Item(16, 18).Value = "Hortaliza"
